$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the three test-case title cells (shared-string renames / reorder) ---
$ws.Range("A1").Value2 = "TCUO1_validateValidEmployeeUpdate"
$ws.Range("A5").Value2 = "TCU02_validateUpdationWithInvalidEmployeeId"
$ws.Range("A9").Value2 = "TCU03_validateDataUpdationOfValidEmployeeWithInvalidData"

# --- Swap the sample-data rows between the "Invalid Employee Id" and
#     "Invalid Data" sections (row 7 <-> row 11) ---
$ws.Range("A7").Value2 = "Manju"
$ws.Range("B7").Value2 = "Dev"
$ws.Range("A11").Value2 = "$%^"
$ws.Range("B11").Value2 = "*&^"

# --- Re-style the section title cells: Consolas 10pt black on yellow fill ---
$titles = "A1","A5","A9"
foreach ($addr in $titles) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Consolas"
    $r.Font.Size = 10
    $r.Font.Color = 0
    $r.Interior.Color = 65535
}

# --- Column layout: column A wider & explicit, rest default ---
# (23.6 is the closest achievable ColumnWidth that rounds to the target
#  stored width of 24.44140625 given this engine's pixel quantization)
$ws.Columns.Item(1).ColumnWidth = 23.6

# --- Selection moved to D16 ---
$ws.Range("D16").Select() | Out-Null
